$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected cells to remain plain text so values such as
# "317.38", "3.420" or "1.88%" are not auto-converted into numbers
# or percentages by Excel.
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "E24", "D25", "E25", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "317.38"
$ws.Range("E2").Value = "1.88%"
$ws.Range("D3").Value = "38.02"
$ws.Range("E3").Value = "1.95%"
$ws.Range("D4").Value = "5.175"
$ws.Range("E4").Value = "1.02%"
$ws.Range("E5").Value = "2.13%"
$ws.Range("D6").Value = "8.526"
$ws.Range("E6").Value = "3.14%"
$ws.Range("D7").Value = "1.936"
$ws.Range("E7").Value = "1.76%"
$ws.Range("D8").Value = "2.954"
$ws.Range("E8").Value = "0.53%"
$ws.Range("D9").Value = "0.9439"
$ws.Range("E9").Value = "2.80%"
$ws.Range("D10").Value = "0.1257"
$ws.Range("E10").Value = "4.89%"
$ws.Range("E11").Value = "1.31%"
$ws.Range("D12").Value = "0.09079"
$ws.Range("E12").Value = "1.32%"
$ws.Range("D13").Value = "0.03415"
$ws.Range("E13").Value = "1.85%"
$ws.Range("D14").Value = "0.09519"
$ws.Range("E14").Value = "-0.74%"
$ws.Range("D15").Value = "0.001359"
$ws.Range("E15").Value = "-1.83%"
$ws.Range("D16").Value = "0.006109"
$ws.Range("E16").Value = "6.21%"
$ws.Range("D17").Value = "3.420"
$ws.Range("E17").Value = "-3.04%"
$ws.Range("E18").Value = "1.18%"
$ws.Range("D19").Value = "0.3516"
$ws.Range("E19").Value = "2.19%"
$ws.Range("D20").Value = "6.584"
$ws.Range("E20").Value = "25.01%"
$ws.Range("D21").Value = "0.1307"
$ws.Range("E21").Value = "1.72%"
$ws.Range("D22").Value = "0.2415"
$ws.Range("E22").Value = "-6.92%"
$ws.Range("D23").Value = "0.04368"
$ws.Range("E23").Value = "0.25%"
$ws.Range("E24").Value = "-2.10%"
$ws.Range("D25").Value = "0.004264"
$ws.Range("E25").Value = "-8.53%"
$ws.Range("E26").Value = "-2.86%"
$ws.Range("D27").Value = "0.0003972"
$ws.Range("E27").Value = "-0.65%"
$ws.Range("D39").Value = "0.02443"
$ws.Range("E39").Value = "8.14%"
$ws.Range("D40").Value = "0.05155"
$ws.Range("E40").Value = "2.32%"
$ws.Range("D41").Value = "0.007489"
$ws.Range("E41").Value = "0.12%"
$ws.Range("D42").Value = "0.1399"
$ws.Range("E42").Value = "3.74%"
$ws.Range("D43").Value = "0.008563"
$ws.Range("E43").Value = "-5.54%"
$ws.Range("D44").Value = "0.002031"
$ws.Range("E44").Value = "3.90%"
$ws.Range("D45").Value = "0.008739"
$ws.Range("E45").Value = "-6.10%"
$ws.Range("D46").Value = "0.00006459"
$ws.Range("E46").Value = "-1.65%"
$ws.Range("D47").Value = "0.00000000747"
$ws.Range("E47").Value = "-0.66%"
$ws.Range("E48").Value = "-14.88%"
$ws.Range("D49").Value = "0.001682"
$ws.Range("E49").Value = "67.84%"
$ws.Range("D50").Value = "0.00002091"
$ws.Range("E50").Value = "-0.66%"
$ws.Range("D51").Value = "0.0001991"
$ws.Range("E51").Value = "-0.66%"
